$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data in columns B:G (rows 2-10) is shifted down by one row (to rows 3-11),
# and a new computed row of values is written into row 2. The previous row 11
# values fall off the bottom of the table.

# Step 1: read old values for rows 2-10, columns B-G (cols 2-7)
$vals = @{}
for ($r = 2; $r -le 10; $r++) {
    for ($c = 2; $c -le 7; $c++) {
        $vals["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# Step 2: write them shifted down by one row (row r -> row r+1)
for ($r = 2; $r -le 10; $r++) {
    for ($c = 2; $c -le 7; $c++) {
        $ws.Cells.Item($r + 1, $c).Value = $vals["$r,$c"]
    }
}

# Step 3: write the new computed values into row 2 (B2:G2)
$ws.Range("B2").Value = 0.01132367786385012
$ws.Range("C2").Value = 2.289151444524298
$ws.Range("D2").Value = 20.74398971997876
$ws.Range("E2").Value = 4.554557027854494
$ws.Range("F2").Value = 4.656905011860751
$ws.Range("G2").Value = 23
